# Enhance Siege Analytics descriptions: add three new bullet points right
# after the "Data Engineering and Infrastructure Architecture" subheading
# paragraph, before the existing "• Architect enterprise-scale ..." bullet.

$d = $word.ActiveDocument

$anchorText = "Data Engineering and Infrastructure Architecture"
$newBullets = @(
    "• Architected data infrastructure processing 15+ billion voter records to support meta-analytical voter file corrections",
    "• Built scalable ETL pipelines enabling analysis of 50,000+ electoral boundaries across all levels of government",
    "• Developed Python boundary estimation algorithm that reduced mapping costs by 75% for 200+ organizations"
)

# Locate the anchor paragraph (the subheading directly under the
# "PARTNER - Siege Analytics" job title).
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r") -eq $anchorText) {
        $anchor = $p
        break
    }
}

if ($anchor -eq $null) {
    Write-Output "ERROR: anchor paragraph not found"
} else {
    $current = $anchor
    foreach ($bulletText in $newBullets) {
        $current.Range.InsertParagraphAfter()
        $current = $current.Next()
        $current.Range.Text = $bulletText
    }
    Write-Output "Inserted $($newBullets.Count) bullet paragraphs after '$anchorText'"
}
